$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 94; this shifts the existing rows 94..182
# down to 95..183 (carrying all their values/styles with them) and
# extends the sheet dimension to A1:R183 automatically.
$ws.Range("A94").EntireRow.Insert()

# Populate the newly inserted row 94 with the new daily record
# (same market/category/quality metadata as its neighbours, new date
# and price figures).
$ws.Range("A94").Value = 9
$ws.Range("B94").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C94").Value = "Metropolitana"
$ws.Range("D94").Value = 44512
$ws.Range("E94").Value = 13
$ws.Range("F94").Value = 300000001
$ws.Range("G94").Value = "Rabanito"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 7900
$ws.Range("K94").Value = 2500
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = 2747
$ws.Range("N94").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O94").Value = "Provincia de Chacabuco"
$ws.Range("P94").Value = 27
$ws.Range("Q94").Value = 100
$ws.Range("R94").Value = "Hortaliza"
